$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JobType")

# Copy formatting from the last data row (A2) down to the new rows A3:A5
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Fill in the new JobType values
$ws.Range("A3").Value = "Strategy"
$ws.Range("A4").Value = "Post Merger Integration"
$ws.Range("A5").Value = "Valuation Advisory"

# Match the saved selection state
$ws.Range("A3:A5").Select()
